# Add "User Account" test-data columns (signin_title / create_title) to the
# customer sample data sheet, update a couple of existing sample values, and
# refresh row heights / selection to match the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the two added columns -----------------------
$ws.Range("N1").Value = "signin_title"
$ws.Range("O1").Value = "create_title"

# --- New data cells for row 2 (vinayms) --------------------------------
$ws.Range("N2").Value = "Sign In with Email"
$ws.Range("O2").Value = "Create your account"

# --- New data cells for row 3 (sanju) -----------------------------------
$ws.Range("N3").Value = "Sign In with Email"
$ws.Range("O3").Value = "Create your account"

# --- Updated sample values ----------------------------------------------
$ws.Range("B2").Value = "vinayms"
$ws.Range("E2").Value = 278439847
$ws.Range("B3").Value = "sanju"

# --- Row heights ----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 25.5
$ws.Rows.Item(2).RowHeight = 14.9
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(4).RowHeight = 12.8

# --- View / selection -----------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$null = $ws.Range("N5").Select()
